# Update the Z (objective), Gap and Time columns for each instance row
# with freshly recomputed solver results (see commit message: "Add the
# models to the output and the data to the excel of each instance").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lib_2")

$ws.Range("F2").Value = 105
$ws.Range("D3").Value = 7913
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 67
$ws.Range("F4").Value = 94
$ws.Range("F5").Value = 354
$ws.Range("F6").Value = 74
$ws.Range("F7").Value = 59
$ws.Range("D8").Value = 9487.999999999998
$ws.Range("E8").Value = 0.000000000001818989403545856072
$ws.Range("F8").Value = 141
$ws.Range("F9").Value = 423
$ws.Range("F10").Value = 76
$ws.Range("F11").Value = 61
$ws.Range("F12").Value = 85
$ws.Range("F13").Value = 131
$ws.Range("D14").Value = 8251.999999999998
$ws.Range("F14").Value = 201
$ws.Range("F15").Value = 101
$ws.Range("F16").Value = 171
$ws.Range("E17").Value = 0.000000000001818989403545856072
$ws.Range("F17").Value = 270
$ws.Range("F18").Value = 169
$ws.Range("D19").Value = 7124.999999999999
$ws.Range("E19").Value = 0.000000000000909494701772928238
$ws.Range("F19").Value = 94
$ws.Range("D20").Value = 8885.999999999993
$ws.Range("E20").Value = 0.000000000007275957614183425903
$ws.Range("F20").Value = 488
$ws.Range("F21").Value = 782
$ws.Range("F22").Value = 115
$ws.Range("F23").Value = 96
$ws.Range("D24").Value = 8746
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 226
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 433
$ws.Range("F26").Value = 4596
$ws.Range("F27").Value = 1769
$ws.Range("D28").Value = 12322
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 13082
$ws.Range("D29").Value = 13722
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 8487
$ws.Range("D30").Value = 12371
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 16619
$ws.Range("D31").Value = 11331
$ws.Range("E31").Value = 0.000000000001818989403545856072
$ws.Range("F31").Value = 33944
$ws.Range("D32").Value = 13331
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 52446
$ws.Range("D33").Value = 15331
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 39277
$ws.Range("F34").Value = 2727
$ws.Range("F35").Value = 1447
$ws.Range("F36").Value = 2205
$ws.Range("F37").Value = 5706
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 1070
$ws.Range("F39").Value = 893
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 1111
$ws.Range("F41").Value = 1488
$ws.Range("F42").Value = 1516
$ws.Range("F43").Value = 2878
$ws.Range("F44").Value = 463
$ws.Range("F45").Value = 356
$ws.Range("F46").Value = 426
$ws.Range("D47").Value = 5650.999999999999
$ws.Range("E47").Value = 0.000000000000909494701772928238
$ws.Range("F47").Value = 815
$ws.Range("F48").Value = 86
$ws.Range("F49").Value = 233
$ws.Range("F50").Value = 446
$ws.Range("D51").Value = 8741
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 3476
$ws.Range("D52").Value = 7414
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 4189
$ws.Range("D53").Value = 9178.000000000002
$ws.Range("E53").Value = 0.000000000001818989403545856072
$ws.Range("F53").Value = 366
$ws.Range("F54").Value = 245
$ws.Range("F55").Value = 150
$ws.Range("D56").Value = 7654.000000000001
$ws.Range("E56").Value = 0.000000000000909494701772928238
$ws.Range("F56").Value = 816
$ws.Range("D57").Value = 21103
$ws.Range("E57").Value = 0
$ws.Range("F57").Value = 12071
$ws.Range("D58").Value = 26039
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 72795
$ws.Range("D59").Value = 37239
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 673029
$ws.Range("D60").Value = 27282
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 22506
$ws.Range("F61").Value = 1341
$ws.Range("D62").Value = 24454.00000000001
$ws.Range("E62").Value = 0.000000000010913936421275140471
$ws.Range("F62").Value = 1385
$ws.Range("D63").Value = 32643
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 25449
$ws.Range("F64").Value = 1932
$ws.Range("F65").Value = 1458
$ws.Range("F66").Value = 1679
$ws.Range("D67").Value = 31415
$ws.Range("E67").Value = 0
$ws.Range("F67").Value = 30557
$ws.Range("D68").Value = 24848
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 4923
$ws.Range("F69").Value = 1058
$ws.Range("F70").Value = 1460
$ws.Range("D71").Value = 32321
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 16312
$ws.Range("F72").Value = 3642
